$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Assets section: shorten the RPAChallenge_Path value
$ws.Range("C9").Value = "Data\Input\challenge.xlsx"

# Clear the fill style previously applied to C10 (remove the highlighted style),
# restoring it to the default cell style while keeping its text value.
$ws.Range("C10").Value = "https://rpachallenge.com/"
$ws.Range("C10").Style = "Normal"

# Bump the package version shown in the deployment sheet
$ws.Range("C17").Value = "1.0.3"

# Update the active selection to C13 as left by the author when saving
$ws.Range("C13").Select()
